$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update LH2 sweep results for the successful 0.5s-timestep run
$ws.Range("B2").Value = 32.597407196161399
$ws.Range("C2").Value = 15.501639725230945
$ws.Range("D2").Value = 0.47554824320679051
$ws.Range("E2").Value = 29.834186891008578
$ws.Range("F2").Value = 14.815249129564203
$ws.Range("G2").Value = 0.49658632171501277
$ws.Range("H2").Value = 280.5
$ws.Range("I2").Value = 246
$ws.Range("B3").Value = 38.352861789346932
$ws.Range("C3").Value = 21.701605628224257
$ws.Range("D3").Value = 0.56584058179074903
$ws.Range("E3").Value = 35.358243361392518
$ws.Range("F3").Value = 19.903971626709009
$ws.Range("G3").Value = 0.56292309047349487
$ws.Range("H3").Value = 286
$ws.Range("I3").Value = 253.5
$ws.Range("B4").Value = 44.084683388867056
$ws.Range("C4").Value = 27.620115003751266
$ws.Range("D4").Value = 0.62652406415436168
$ws.Range("E4").Value = 41.409843570095433
$ws.Range("F4").Value = 25.557634818811671
$ws.Range("G4").Value = 0.61718742732146892
$ws.Range("H4").Value = 291.5
$ws.Range("I4").Value = 260.5
$ws.Range("B5").Value = 49.842921094285622
$ws.Range("C5").Value = 33.295805580683037
$ws.Range("D5").Value = 0.66801473207597228
$ws.Range("E5").Value = 47.035552786530197
$ws.Range("F5").Value = 30.862774596646975
$ws.Range("G5").Value = 0.65615843267999818
$ws.Range("H5").Value = 297
$ws.Range("I5").Value = 267
$ws.Range("B6").Value = 55.510466207112884
$ws.Range("C6").Value = 39.007391529528959
$ws.Range("D6").Value = 0.70270336739724071
$ws.Range("E6").Value = 53.186949049256185
$ws.Range("F6").Value = 36.639759166959827
$ws.Range("G6").Value = 0.68888627420662762
$ws.Range("H6").Value = 302
$ws.Range("I6").Value = 273.5
$ws.Range("B7").Value = 61.458172021283971
$ws.Range("C7").Value = 44.825746118143144
$ws.Range("D7").Value = 0.7293699868362381
$ws.Range("E7").Value = 58.831895673601551
$ws.Range("F7").Value = 42.016000492717481
$ws.Range("G7").Value = 0.7141704344497346
$ws.Range("H7").Value = 307.5
$ws.Range("I7").Value = 279.5
$ws.Range("B8").Value = 67.508528696510226
$ws.Range("C8").Value = 50.824996513473721
$ws.Range("D8").Value = 0.7528677856684064
$ws.Range("E8").Value = 64.731618236745291
$ws.Range("F8").Value = 47.62839112998536
$ws.Range("G8").Value = 0.73578248817745173
$ws.Range("H8").Value = 313.5
$ws.Range("I8").Value = 285.5
$ws.Range("B9").Value = 73.300529854644552
$ws.Range("C9").Value = 56.489836683650964
$ws.Range("D9").Value = 0.77066068684183719
$ws.Range("E9").Value = 70.662341844340901
$ws.Range("F9").Value = 53.163598152276499
$ws.Range("G9").Value = 0.7523611129303972
$ws.Range("H9").Value = 317
$ws.Range("I9").Value = 291.5
$ws.Range("B10").Value = 79.400746093425212
$ws.Range("C10").Value = 62.637738852153305
$ws.Range("D10").Value = 0.7888809858090241
$ws.Range("E10").Value = 76.387105122079021
$ws.Range("F10").Value = 58.515450115818716
$ws.Range("G10").Value = 0.76603832573968478
$ws.Range("H10").Value = 323
$ws.Range("I10").Value = 296.5
$ws.Range("B11").Value = 85.486097537657372
$ws.Range("C11").Value = 68.511372690086148
$ws.Range("D11").Value = 0.80143291907676906
$ws.Range("E11").Value = 82.618366810570834
$ws.Range("F11").Value = 64.4198641185978
$ws.Range("G11").Value = 0.77972812348495157
$ws.Range("H11").Value = 326.5
$ws.Range("I11").Value = 302

# Re-fit the columns to the refreshed values
$ws.Columns.Item(1).ColumnWidth = 17.583781
$ws.Columns.Item(2).ColumnWidth = 32.417156
$ws.Columns.Item(3).ColumnWidth = 30.917469
$ws.Columns.Item(4).ColumnWidth = 37.750938
$ws.Columns.Item(5).ColumnWidth = 31.084
$ws.Columns.Item(6).ColumnWidth = 29.584312
$ws.Columns.Item(7).ColumnWidth = 36.583781
$ws.Columns.Item(8).ColumnWidth = 31.417156
$ws.Columns.Item(9).ColumnWidth = 30.084
